{"js": "// Apply three small wording fixes around \"displej\" (display) in the document.\n// 1) \"...pracovat s displeji\" -> \"...pracovat s displejem\"\n// 2) \"vypisovat hodnoty a data na displej\" -> \"vypisovat hodnoty na displeji\"\n// 3) \"zprovoznit displeje, \" -> \"zprovoznit displej, \"\n\nconst body = context.document.body;\n\n// --- Change 1: displeji -> displejem ------------------------------------\n// Anchor on the preceding words (with the non-breaking space that's really\n// in the document) so this only ever targets this one bullet, never the\n// \"displeji\" substring that later shows up inside \"na displeji\".\nconst results1 = body.search(\"pracovat s\\u00A0displeji\", { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\nif (results1.items.length > 0) {\n  results1.items[0].insertText(\"pracovat s\\u00A0displejem\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: vypisovat hodnoty a data na displej -> vypisovat hodnoty na displeji\nconst results2 = body.search(\"vypisovat hodnoty a data na displej\", { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\nif (results2.items.length > 0) {\n  results2.items[0].insertText(\"vypisovat hodnoty na displeji\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 3: zprovoznit displeje,  -> zprovoznit displej,  ------------\nconst results3 = body.search(\"zprovoznit displeje, \", { matchCase: true });\nresults3.load(\"items\");\nawait context.sync();\nif (results3.items.length > 0) {\n  results3.items[0].insertText(\"zprovoznit displej, \", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply three small wording fixes around \"displej\" (display) in the document.\n# 1) \"...pracovat s displeji\" -> \"...pracovat s displejem\"\n# 2) \"vypisovat hodnoty a data na displej\" -> \"vypisovat hodnoty na displeji\"\n# 3) \"zprovoznit displeje, \" -> \"zprovoznit displej, \"\n\n$d = $word.ActiveDocument\n\nfunction Replace-OneMatch($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0              # wdFindStop - do not wrap around the document\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    # wdReplaceOne = 1 -> only the (single, targeted) match gets replaced\n    $find.Execute($null, $true, $null, $null, $null, $null, $true, 1, $false, $replaceText, 1)\n}\n\n# Anchor on the preceding words (including the non-breaking space that is\n# really present in the document) so this only ever targets this one\n# bullet, never the \"displeji\" substring that later shows up inside\n# \"na displeji\" once change 2 below has run.\n$nbsp = [char]0x00A0\n$search1 = \"pracovat s\" + $nbsp + \"displeji\"\n$replace1 = \"pracovat s\" + $nbsp + \"displejem\"\nReplace-OneMatch $search1 $replace1\n\nReplace-OneMatch \"vypisovat hodnoty a data na displej\" \"vypisovat hodnoty na displeji\"\n\nReplace-OneMatch \"zprovoznit displeje, \" \"zprovoznit displej, \"\n"}
